# "Having issue reading in county data"
# Rename the two sheets to reflect what they actually hold, and leave the
# selection on the county sheet parked where the author was last looking
# while debugging the read issue.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "state"
$wb.Worksheets.Item(2).Name = "county"

$ws = $wb.Worksheets.Item("county")
$ws.Activate()
$ws.Range("K38").Select()
